# Refresh the "Price" (column D) and "Volume(1h)" (column E) figures on
# the cryptos sheet with the values from the latest scheduled pull.
#
# Column D values frequently look like plain numbers (e.g. "523.07" or
# "59.006.03", the latter being a thousands-dotted price, not a number).
# Assigning such a string straight to .Value lets Excel reinterpret it as
# a float (e.g. 523.07 becomes 523.0700000000005, and "59.006.03" would
# fail to round-trip at all), so for column D we briefly force the cell
# to Text format, assign the literal string, then restore the Normal
# style so the cell's formatting is left exactly as it was before.
# Column E values (e.g. "  +1.71%  ") always contain spaces/a percent
# sign and are never mistaken for numbers, so they can be set directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.006.03'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.71%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.587.61'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '523.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.09'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.25%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.564'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.598.83'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.54'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.02%  '
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('E12').Value = '  +1.14%  '
$ws.Range('E13').Value = '  +3.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.044.47'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '58.941.04'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.48'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.584.73'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.52%  '
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '338.33'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.48%  '
$ws.Range('E20').Value = '  -0.23%  '
$ws.Range('E21').Value = '  -2.03%  '
$ws.Range('E22').Value = '  +2.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.13'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('E25').Value = '  +1.40%  '
$ws.Range('E26').Value = '  +0.32%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.01'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0724'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.78%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.90'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.54%  '
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.68'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '149.06'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('E35').Value = '  -0.77%  '
$ws.Range('E36').Value = '  -1.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '36.79'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.23%  '
$ws.Range('E38').Value = '  +0.90%  '
$ws.Range('E39').Value = '  -0.81%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.814'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.04%  '
$ws.Range('E41').Value = '  -0.53%  '
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '271.74'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('E44').Value = '  +1.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0955'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.33%  '
$ws.Range('E46').Value = '  +0.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0516'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.41'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.965.06'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('E50').Value = '  -2.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0220'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.39%  '
